$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1966205837173579
$ws.Range("C2").Value = 0.5345622119815668
$ws.Range("J2").Value = 0.01228878648233487
$ws.Range("P2").Value = 0.1490015360983103
$ws.Range("S2").Value = 0.1075268817204301
$ws.Range("B3").Value = 0.005730659025787965
$ws.Range("C3").Value = 0.002865329512893983
$ws.Range("J3").Value = 0.01146131805157593
$ws.Range("P3").Value = 0.7335243553008596
$ws.Range("S3").Value = 0.2464183381088825
$ws.Range("J4").Value = 0.04901960784313725
$ws.Range("P4").Value = 0.696078431372549
$ws.Range("S4").Value = 0.2549019607843137
$ws.Range("B6").Value = 0.07028112449799197
$ws.Range("D6").Value = 0.02008032128514056
$ws.Range("F6").Value = 0.07429718875502007
$ws.Range("J6").Value = 0.2931726907630522
$ws.Range("O6").Value = 0.02208835341365462
$ws.Range("Q6").Value = 0.1365461847389558
$ws.Range("R6").Value = 0.06224899598393574
$ws.Range("S6").Value = 0.321285140562249
$ws.Range("B7").Value = 0.1253481894150418
$ws.Range("D7").Value = 0.02506963788300836
$ws.Range("F7").Value = 0.03342618384401114
$ws.Range("J7").Value = 0.1337047353760446
$ws.Range("O7").Value = 0.03064066852367688
$ws.Range("Q7").Value = 0.1587743732590529
$ws.Range("R7").Value = 0.07242339832869081
$ws.Range("S7").Value = 0.4206128133704735
$ws.Range("B8").Value = 0.08884297520661157
$ws.Range("D8").Value = 0.02789256198347108
$ws.Range("E8").Value = 0.001033057851239669
$ws.Range("F8").Value = 0.0640495867768595
$ws.Range("J8").Value = 0.1022727272727273
$ws.Range("O8").Value = 0.01756198347107438
$ws.Range("Q8").Value = 0.1652892561983471
$ws.Range("R8").Value = 0.1074380165289256
$ws.Range("S8").Value = 0.4256198347107438
$ws.Range("B9").Value = 0.08571428571428572
$ws.Range("D9").Value = 0.01714285714285714
$ws.Range("F9").Value = 0.04952380952380953
$ws.Range("J9").Value = 0.1142857142857143
$ws.Range("O9").Value = 0.02666666666666667
$ws.Range("Q9").Value = 0.1961904761904762
$ws.Range("R9").Value = 0.09904761904761905
$ws.Range("S9").Value = 0.4114285714285714
$ws.Range("B10").Value = 0.109599395313681
$ws.Range("D10").Value = 0.0200302343159486
$ws.Range("E10").Value = 0.0003779289493575208
$ws.Range("F10").Value = 0.06349206349206349
$ws.Range("J10").Value = 0.1250944822373394
$ws.Range("O10").Value = 0.01473922902494331
$ws.Range("Q10").Value = 0.1806500377928949
$ws.Range("R10").Value = 0.09108087679516251
$ws.Range("S10").Value = 0.3949357520786092
$ws.Range("G11").Value = 0.1474245115452931
$ws.Range("J11").Value = 0.1119005328596803
$ws.Range("K11").Value = 0.2042628774422735
$ws.Range("L11").Value = 0.5150976909413855
$ws.Range("S11").Value = 0.02131438721136767
$ws.Range("G12").Value = 0.7375415282392026
$ws.Range("J12").Value = 0.1827242524916944
$ws.Range("K12").Value = 0.0132890365448505
$ws.Range("L12").Value = 0.02990033222591362
$ws.Range("S12").Value = 0.03654485049833887
$ws.Range("F15").Value = 0.03155818540433925
$ws.Range("H15").Value = 0.1538461538461539
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3510848126232742
$ws.Range("K15").Value = 0.04536489151873768
$ws.Range("M15").Value = 0.01380670611439842
$ws.Range("O15").Value = 0.07495069033530571
$ws.Range("S15").Value = 0.252465483234714
$ws.Range("F16").Value = 0.02926829268292683
$ws.Range("H16").Value = 0.1731707317073171
$ws.Range("I16").Value = 0.08292682926829269
$ws.Range("J16").Value = 0.3804878048780488
$ws.Range("K16").Value = 0.0975609756097561
$ws.Range("M16").Value = 0.02195121951219512
$ws.Range("N16").Value = 0.004878048780487805
$ws.Range("O16").Value = 0.06829268292682927
$ws.Range("S16").Value = 0.1414634146341463
$ws.Range("F17").Value = 0.02076124567474048
$ws.Range("H17").Value = 0.1949250288350634
$ws.Range("I17").Value = 0.1003460207612457
$ws.Range("J17").Value = 0.3956170703575548
$ws.Range("K17").Value = 0.08650519031141868
$ws.Range("M17").Value = 0.01614763552479815
$ws.Range("N17").Value = 0.001153402537485583
$ws.Range("O17").Value = 0.06574394463667819
$ws.Range("S17").Value = 0.118800461361015
$ws.Range("F18").Value = 0.04424778761061947
$ws.Range("H18").Value = 0.1703539823008849
$ws.Range("I18").Value = 0.1283185840707965
$ws.Range("J18").Value = 0.336283185840708
$ws.Range("K18").Value = 0.09292035398230089
$ws.Range("M18").Value = 0.008849557522123894
$ws.Range("N18").Value = 0.002212389380530973
$ws.Range("O18").Value = 0.09070796460176991
$ws.Range("S18").Value = 0.1261061946902655
$ws.Range("F19").Value = 0.02485795454545454
$ws.Range("H19").Value = 0.2077414772727273
$ws.Range("I19").Value = 0.109375
$ws.Range("J19").Value = 0.3597301136363636
$ws.Range("K19").Value = 0.09019886363636363
$ws.Range("M19").Value = 0.02024147727272727
$ws.Range("N19").Value = 0.0007102272727272727
$ws.Range("O19").Value = 0.06924715909090909
$ws.Range("S19").Value = 0.11789772727272728
